$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$halfPi = 1.5707963267948966

for ($row = 1; $row -le 17; $row++) {
    $bVal = $ws.Cells.Item($row, 2).Value()
    $cVal = $ws.Cells.Item($row, 3).Value()
    $dVal = $ws.Cells.Item($row, 4).Value()
    $eVal = $ws.Cells.Item($row, 5).Value()
    $fVal = $ws.Cells.Item($row, 6).Value()

    $ws.Cells.Item($row, 2).Value = -$bVal
    $ws.Cells.Item($row, 3).Value = $halfPi - $cVal
    $ws.Cells.Item($row, 4).Value = -$dVal
    $ws.Cells.Item($row, 5).Value = $halfPi - $eVal
    $ws.Cells.Item($row, 6).Value = -$fVal
}
